$d = $word.ActiveDocument
$d.Content.Find.Execute("855×6=5130", $false, $false, $false, $false, $false, $true, 1, $false, "399×3=1197", 2) | Out-Null
$d.Content.Find.Execute("545×6=3270", $false, $false, $false, $false, $false, $true, 1, $false, "198×5=990", 2) | Out-Null
$d.Content.Find.Execute("214×8=1712", $false, $false, $false, $false, $false, $true, 1, $false, "716×5=3580", 2) | Out-Null
$d.Content.Find.Execute("350×4=1400", $false, $false, $false, $false, $false, $true, 1, $false, "332×5=1660", 2) | Out-Null
$d.Content.Find.Execute("163×7=1141", $false, $false, $false, $false, $false, $true, 1, $false, "981×5=4905", 2) | Out-Null
$d.Content.Find.Execute("483×9=4347", $false, $false, $false, $false, $false, $true, 1, $false, "551×2=1102", 2) | Out-Null
$d.Content.Find.Execute("361×7=2527", $false, $false, $false, $false, $false, $true, 1, $false, "802×8=6416", 2) | Out-Null
$d.Content.Find.Execute("274×6=1644", $false, $false, $false, $false, $false, $true, 1, $false, "588×3=1764", 2) | Out-Null
$d.Content.Find.Execute("625×7=4375", $false, $false, $false, $false, $false, $true, 1, $false, "977×8=7816", 2) | Out-Null
$d.Content.Find.Execute("309×3=927", $false, $false, $false, $false, $false, $true, 1, $false, "215×3=645", 2) | Out-Null
$d.Content.Find.Execute("988×2=1976", $false, $false, $false, $false, $false, $true, 1, $false, "380×5=1900", 2) | Out-Null
$d.Content.Find.Execute("710×9=6390", $false, $false, $false, $false, $false, $true, 1, $false, "520×7=3640", 2) | Out-Null
$d.Content.Find.Execute("502×8=4016", $false, $false, $false, $false, $false, $true, 1, $false, "384×8=3072", 2) | Out-Null
$d.Content.Find.Execute("809×6=4854", $false, $false, $false, $false, $false, $true, 1, $false, "125×9=1125", 2) | Out-Null
$d.Content.Find.Execute("900×2=1800", $false, $false, $false, $false, $false, $true, 1, $false, "200×8=1600", 2) | Out-Null
$d.Content.Find.Execute("808×3=2424", $false, $false, $false, $false, $false, $true, 1, $false, "675×7=4725", 2) | Out-Null
$d.Content.Find.Execute("892×7=6244", $false, $false, $false, $false, $false, $true, 1, $false, "825×2=1650", 2) | Out-Null
$d.Content.Find.Execute("780×4=3120", $false, $false, $false, $false, $false, $true, 1, $false, "600×6=3600", 2) | Out-Null
$d.Content.Find.Execute("310×5=1550", $false, $false, $false, $false, $false, $true, 1, $false, "818×3=2454", 2) | Out-Null
$d.Content.Find.Execute("272×5=1360", $false, $false, $false, $false, $false, $true, 1, $false, "633×6=3798", 2) | Out-Null
$d.Content.Find.Execute("937×8=7496", $false, $false, $false, $false, $false, $true, 1, $false, "863×9=7767", 2) | Out-Null
$d.Content.Find.Execute("681×9=6129", $false, $false, $false, $false, $false, $true, 1, $false, "783×5=3915", 2) | Out-Null
$d.Content.Find.Execute("636×6=3816", $false, $false, $false, $false, $false, $true, 1, $false, "293×4=1172", 2) | Out-Null
$d.Content.Find.Execute("930×5=4650", $false, $false, $false, $false, $false, $true, 1, $false, "311×3=933", 2) | Out-Null
$d.Content.Find.Execute("572×4=2288", $false, $false, $false, $false, $false, $true, 1, $false, "301×7=2107", 2) | Out-Null
